$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StatQuery" column (C) for rows 2-5 contains a Cypher query that
# filters on the study designation 'COT007B'. The study code was corrected
# to 'COTC007B'. Update the text in C2:C5 accordingly (other columns are
# untouched; only the shared-string table gets reshuffled by Excel itself).

$newQuery = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE s.clinical_study_designation IN ['COTC007B'] and demo.breed in ['Vizsla'] and demo.sex in['Male']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery
$ws.Range("C5").Value = $newQuery

# Restore the selection to the cell that was active when the author saved
# the file.
$ws.Range("D5").Select()
